# Auto-generated edit script
# Applies numeric value updates to the Balance.xlsx 'BALANCE' workbook sheets
# produced by a re-run of the ATR-72 'ATsi' sizing/balance case (test #01, Sandbox 2).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C3").Value = 12.141605723464071
$ws.Range("C5").Value = 0.1170958407255597
$ws.Range("C7").Value = 48.8697450578075
$ws.Range("C9").Value = 5.046722739735554
$ws.Range("C13").Value = 11.785103940501518
$ws.Range("C15").Value = 0.265163019047306
$ws.Range("C17").Value = 33.50484700775513
$ws.Range("C19").Value = 11.428281565519924
$ws.Range("C23").Value = 11.785103940501518
$ws.Range("C25").Value = 0.265163019047306
$ws.Range("C27").Value = 33.50484700775513
$ws.Range("C29").Value = 11.428281565519924
$ws.Range("C33").Value = 11.785103940501518
$ws.Range("C35").Value = 0.265163019047306
$ws.Range("C37").Value = 33.50484700775513
$ws.Range("C39").Value = 11.428281565519924
$ws.Range("C43").Value = 11.97639295538773
$ws.Range("C45").Value = 0.1802658333969227
$ws.Range("C47").Value = 41.74922712669088
$ws.Range("C49").Value = 7.769291163243255
$ws.Range("C53").Value = 11.963789830651873
$ws.Range("C55").Value = 0.46880783296563133
$ws.Range("C57").Value = 41.20604407422656
$ws.Range("C59").Value = 20.20518522719279
$ws.Range("C62").Value = 18.778612969176336
$ws.Range("C63").Value = 41.20604407422656
$ws.Range("C64").Value = 54.32676647955055
$ws.Range("C69").Value = 47870.72870372752
$ws.Range("C70").Value = 745048.6735852612
$ws.Range("C71").Value = 697177.9448815337
$ws.Range("C76").Value = 48532.61138617322

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C5").Value = 1.7164117478640433
$ws.Range("C6").Value = 1.7164117478640046

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C5").Value = 33.26255957187139
$ws.Range("C6").Value = 33.26255957187139
$ws.Range("C7").Value = 0.9794823019718901
$ws.Range("C8").Value = 0.9794823019718901
$ws.Range("C9").Value = 11.779482301971889
$ws.Range("C10").Value = 11.779482301971889
$ws.Range("C12").Value = 5.409999999999998
$ws.Range("C13").Value = 5.4099999999999975
$ws.Range("C15").Value = 68.95852434675184
$ws.Range("C16").Value = 68.95852434675184
$ws.Range("C23").Value = 1.104420538817831
$ws.Range("C24").Value = 0.8545440651259495
$ws.Range("C27").Value = 5.409999999999998
$ws.Range("C28").Value = 10.998778735632179

$ws = $wb.Worksheets.Item("FUEL TANK")
$ws.Range("C5").Value = 39.07656408472137
$ws.Range("C6").Value = 39.07656408472133
$ws.Range("C7").Value = 1.1143808849778667
$ws.Range("C8").Value = 1.1143808849778667
$ws.Range("C9").Value = 11.914380884977867
$ws.Range("C10").Value = 11.914380884977865
$ws.Range("C15").Value = 68.95852434675184
$ws.Range("C16").Value = 68.95852434675184

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C5").Value = 654.0912358990008
$ws.Range("C6").Value = 654.0912358990008
$ws.Range("C15").Value = 247.27664849190882
$ws.Range("C16").Value = 247.27664849190882

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C5").Value = 584.4777395215115
$ws.Range("C6").Value = 584.4777395215115
$ws.Range("C15").Value = 56.028801031735874
$ws.Range("C16").Value = 56.028801031735874

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C5").Value = 6.347948306468693
$ws.Range("C6").Value = 6.347948306468693
$ws.Range("C15").Value = 40.6079510246935
$ws.Range("C16").Value = 40.6079510246935

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C5").Value = -36.75112941025122
$ws.Range("C6").Value = -36.75112941025122
$ws.Range("C15").Value = 40.6079510246935
$ws.Range("C16").Value = 40.6079510246935

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = -19.934136097328047
$ws.Range("C6").Value = -19.934136097328086
$ws.Range("C7").Value = 10.545193809330613
$ws.Range("C8").Value = 10.545193809330613
$ws.Range("C9").Value = 10.545193809330613
$ws.Range("C10").Value = 10.545193809330613
$ws.Range("C15").Value = -186.09319776525322
$ws.Range("C16").Value = -186.09319776525317
$ws.Range("C17").Value = -2.5397999999999996
$ws.Range("C18").Value = -2.539799999999999
$ws.Range("C19").Value = -4.317799999999999
$ws.Range("C20").Value = -4.317799999999998
$ws.Range("C23").Value = 10.545193809330613
$ws.Range("C26").Value = -2.5397999999999996
